$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow the data columns (A:I) to reflect updated content widths
$ws.Columns.Item(1).ColumnWidth = 20.8333333333
$ws.Columns.Item(2).ColumnWidth = 29.5
$ws.Columns.Item(3).ColumnWidth = 28
$ws.Columns.Item(4).ColumnWidth = 34.3333333333
$ws.Columns.Item(5).ColumnWidth = 28.1666666667
$ws.Columns.Item(6).ColumnWidth = 26.8333333333
$ws.Columns.Item(7).ColumnWidth = 33.1666666667
$ws.Columns.Item(8).ColumnWidth = 28.3333333333
$ws.Columns.Item(9).ColumnWidth = 27.1666666667

# Update the feed-temperature sweep results for the new tank-pressure / feed-pressure runs
$ws.Range("A2").Value = 0.01
$ws.Range("B2").Value = 32.447837977294029
$ws.Range("C2").Value = 15.983216035754685
$ws.Range("D2").Value = 0.49258184927264598
$ws.Range("E2").Value = 31.392631591000779
$ws.Range("F2").Value = 14.751070143845119
$ws.Range("G2").Value = 0.46988956950247396
$ws.Range("H2").Value = 510.10000000000002
$ws.Range("I2").Value = 479.69999999999993

$ws.Range("A3").Value = 0.055000000000000007
$ws.Range("B3").Value = 33.066455269485182
$ws.Range("C3").Value = 16.156579950853445
$ws.Range("D3").Value = 0.48860937222270906
$ws.Range("E3").Value = 29.433115708657851
$ws.Range("F3").Value = 14.019504130634449
$ws.Range("G3").Value = 0.47631736542626929
$ws.Range("H3").Value = 271.30000000000001
$ws.Range("I3").Value = 240.10000000000002

$ws.Range("A4").Value = 0.10000000000000001
$ws.Range("B4").Value = 33.428327833143939
$ws.Range("C4").Value = 16.04992266498996
$ws.Range("D4").Value = 0.48012939041110458
$ws.Range("E4").Value = 29.85476044020924
$ws.Range("F4").Value = 14.995006254825931
$ws.Range("G4").Value = 0.50226516755533002
$ws.Range("H4").Value = 232.80000000000001
$ws.Range("I4").Value = 203.59999999999997
